# TextBox가 OwnerObject의 Theater 기준으로 생성되도록 수정
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to reflect its actual purpose.
$ws.Name = "EditorMenuTable"

# Row 2: the focus object used for the TestCharacter theater becomes the
# dedicated "Test_Character_Face" focus object.
$ws.Range("G2").Value = "Test_Character_Face"

# Row 3: fill in the missing NextTheaterAlias (chains to a new theater),
# switch it to play once ("None" instead of "Loop"), and give it an
# OwnerObject/FocusObject pair (PlayerCharacter / Character_Face) so the
# TextBox is anchored off the owner object's theater.
$ws.Range("B3").Value = "TheaterTest3"
$ws.Range("E3").Value = "None"
$ws.Range("F3").Value = "PlayerCharacter"
$ws.Range("G3").Value = "Character_Face"

# Row 4: brand new theater row chained from TheaterTest3.
$ws.Range("A4").Value = "TheaterTest3"
$ws.Range("C4").Value = "<base>테스트으으!3</>"
$ws.Range("D4").Value = "SizeS"
$ws.Range("E4").Value = "Loop"
$ws.Range("F4").Value = "TestCharacter"
$ws.Range("G4").Value = "Test_Character_Face"

# Carry over the existing column formatting (OwnerObject/FocusObject use the
# bold-ish font style; StringData uses vertical-centered alignment) onto the
# newly populated cells.
$ws.Range("F1").Copy() | Out-Null
$ws.Range("F3:G4").PasteSpecial(-4122) | Out-Null
$ws.Range("C2").Copy() | Out-Null
$ws.Range("C4").PasteSpecial(-4122) | Out-Null

# Match the author's final cursor position.
$ws.Range("F4").Select() | Out-Null
